$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A (Região) - rows shift down by one (new region "Amapá" inserted)
$ws.Range("A3").Value = "Amapá"
$ws.Range("A4").Value = "Bahia"
$ws.Range("A5").Value = "Distrito Federal"
$ws.Range("A6").Value = "Alagoas"

# Update column C (Trimestre) - all rows get the new quarter date (stored as plain text)
$cRange = $ws.Range("C2:C9")
$cRange.NumberFormat = "@"
$cRange.Value = "01/07/2025"

# Update column D (Valor) - new numeric values
$ws.Range("D2").Value = 10
$ws.Range("D3").Value = 8.699999999999999
$ws.Range("D4").Value = 8.5
$ws.Range("D5").Value = 8
$ws.Range("D6").Value = 7.7
$ws.Range("D7").Value = 7.7
$ws.Range("D8").Value = 5.6
$ws.Range("D9").Value = 7.8
